$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "new input file order": the old column B ("building category") is removed,
# shifting every later column one to the left.
$ws.Columns("B:B").Delete()

# A new column is inserted right before "cooling system" (now column AD) to
# hold the new "heat pump efficiency" field.
$ws.Columns("AC:AC").Insert()

# Row 1 (headers) / row 2 (units) for the new column.
$ws.Range("AC1").Value = "heat pump efficiency"
$ws.Range("AC2").Value = "[-]"

# "cold emission" -> "cold emission system" (same column, just renamed).
$ws.Range("AE1").Value = "cold emission system"

# Data row (row 3) updates.
$ws.Range("F3").Value = 0
$ws.Range("AB3").Value = "radiator"
$ws.Range("AC3").Value = 0.55000000000000004
$ws.Range("AE3").Value = "air"
$ws.Range("AF3").Value = "105 30"
$ws.Range("AJ3").Value = "180 120"
$ws.Range("AK3").Value = "26 25"

# Restore the view/selection recorded in the saved workbook.
$ws.Range("F4").Select()
$excel.ActiveWindow.ScrollColumn = 4
